# Insert four new list paragraphs right after the paragraph that ends with
# "...URNResource / Contexts / Roles." (the "Switch actions:" bullet), and
# right before the existing blank bullet that precedes "Composition of Case
# classes...". The new paragraphs are:
#   1. an empty bullet
#   2. "(Dimension, Unit, Measure, Value);"
#   3. "Equivalence: Same Distance Measures."
#   4. "Entailment: Dimension, Unit, Measure Values entails other Measures
#       / Values (Time, Speed, Distance)."

$d = $word.ActiveDocument

$anchorText = "Switch actions: Populate Models (RDFS, OWL, Sets, FCA Contexts, Functional MVC / DCI DOM / Others: SaILs). URNResource / Contexts / Roles."

$find = $d.Content
$found = $find.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

# Build an insertion point *after* the matched range's end, using a fresh
# Range object (re-using the Find range directly confuses InsertXML, which
# otherwise clobbers the paragraph that was just matched).
$insertPos = $find.End
$ins = $d.Range($insertPos, $insertPos)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$para1 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:ind w:left=`"600`" w:hanging=`"360`"/></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr></w:r></w:p>"

$para2 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:ind w:left=`"600`" w:hanging=`"360`"/></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">(Dimension, Unit, Measure, Value);</w:t></w:r></w:p>"

$para3 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:ind w:left=`"600`" w:hanging=`"360`"/><w:rPr><w:u w:val=`"none`"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">Equivalence: Same Distance Measures.</w:t></w:r></w:p>"

$para4 = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"3`"/></w:numPr><w:ind w:left=`"600`" w:hanging=`"360`"/><w:rPr><w:u w:val=`"none`"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val=`"0`"/></w:rPr><w:t xml:space=`"preserve`">Entailment: Dimension, Unit, Measure Values entails other Measures / Values (Time, Speed, Distance).</w:t></w:r></w:p>"

$xml = $para1 + $para2 + $para3 + $para4

$ins.InsertXML($xml)

Write-Output "Inserted 4 paragraphs after anchor."
